$wb = $excel.ActiveWorkbook

# Rename the first sheet from "rel.vlažnost" to "rel.vlaznost"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "rel.vlaznost"

# MUV sheet: update selection / active cell
$ws2 = $wb.Worksheets.Item("MUV")
$ws2.Activate()
$ws2.Range("M16").Select()

# temp. sheet: keep its own selection (unchanged: B2), but it should no
# longer be the tab-selected sheet
$ws5 = $wb.Worksheets.Item("temp.")
$ws5.Activate()
$ws5.Range("B2").Select()

# sr.brz. sheet: update selection / active cell
$ws8 = $wb.Worksheets.Item("sr.brz.")
$ws8.Activate()
$ws8.Range("C16").Select()

# rel.vlaznost (first sheet) becomes the active / tab-selected sheet,
# scrolled down and with J21 as the active cell
$ws1.Activate()
$ws1.Range("J21").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
